$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure B2 keeps its text type (value "1" -> "2"), since it looks numeric
$ws.Range("B2").NumberFormat = "@"

$ws.Range("A2").Value = "Netherlands"
$ws.Range("B2").Value = "2"
$ws.Range("C2").Value = "Power"
$ws.Range("D2").Value = -0.424
$ws.Range("G2").Value = -0.1292156119742327
$ws.Range("H2").Value = -0.1660098522167488
$ws.Range("I2").Value = -0.3535460425462049
$ws.Range("J2").Value = -0.3535460425462049
$ws.Range("K2").Value = -17.6
$ws.Range("L2").Value = -0.6669192876089428
$ws.Range("M2").Value = 0.173
$ws.Range("N2").Value = 0.0002037763407423113
$ws.Range("O2").Value = -0.009829545454545454
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = -0
$ws.Range("S2").Value = 0.173
$ws.Range("T2").Value = 1
$ws.Range("U2").Value = 15.474
$ws.Range("V2").Value = 0.01822679246616488
$ws.Range("W2").Value = 2.451243201243201
$ws.Range("X2").Value = 0.04079135508737444
$ws.Range("Y2").Value = 2.410451846155827
$ws.Range("Z2").Value = 0.6487312938813432
$ws.Range("AA2").Value = -0.1433236232278912
$ws.Range("AB2").Value = 0.03333054557380447
$ws.Range("AC2").Value = -0.1766541688016956
$ws.Range("AD2").Value = 73.45
$ws.Range("AE2").Value = 0.1354003139717425
$ws.Range("AF2").Value = 73.58540031397175
$ws.Range("AG2").Value = 58.11140031397174
$ws.Range("AH2").Value = 0.07976258150884874
$ws.Range("AI2").Value = 1.109794978471248
$ws.Range("AJ2").Value = 0.06406415156771753
$ws.Range("AK2").Value = 1.143218560870513
$ws.Range("AL2").Value = 4.275
$ws.Range("AM2").Value = 4.141
$ws.Range("AN2").Value = -15.1349680609932
$ws.Range("AO2").Value = -2.208421052631579
$ws.Range("AP2").Value = -11.97432522439146
$ws.Range("AQ2").Value = -2.279884085969572
$ws.Range("A3").Value = "Netherlands"
$ws.Range("B3").Value = "DGB Group N.V. (ENXTAM:DGB)"
$ws.Range("C3").Value = "Power"
$ws.Range("D3").Value = -0.424
$ws.Range("G3").Value = 0.09470899470899471
$ws.Range("H3").Value = 0.04333333333333333
$ws.Range("I3").Value = -0.006882543004991985
$ws.Range("J3").Value = -0.006882543004991985
$ws.Range("K3").Value = -3.8
$ws.Range("L3").Value = -0.2010582010582011
$ws.Range("M3").Value = 0.173
$ws.Range("N3").Value = 0.02518195050946142
$ws.Range("O3").Value = -0.04552631578947369
$ws.Range("P3").Value = -0
$ws.Range("Q3").Value = -0
$ws.Range("R3").Value = 0
$ws.Range("S3").Value = 0.173
$ws.Range("T3").Value = 1
$ws.Range("U3").Value = 0.274
$ws.Range("V3").Value = 0.03988355167394469
$ws.Range("W3").Value = -0.3247863247863248
$ws.Range("X3").Value = 0.04673970568957293
$ws.Range("Y3").Value = -0.3715260304758977
$ws.Range("Z3").Value = 3.192891002047921
$ws.Range("AA3").Value = -0.02197520963184676
$ws.Range("AB3").Value = 0.03298354988144742
$ws.Range("AC3").Value = -0.05495875951329419
$ws.Range("AD3").Value = 4.95
$ws.Range("AE3").Value = 0.1354003139717425
$ws.Range("AF3").Value = 5.085400313971743
$ws.Range("AG3").Value = 4.811400313971743
$ws.Range("AH3").Value = 0.4253642856298724
$ws.Range("AI3").Value = 0.4066563393648459
$ws.Range("AJ3").Value = 0.4118855774694223
$ws.Range("AK3").Value = 0.3933646345035207
$ws.Range("AL3").Value = 0.275
$ws.Range("AM3").Value = 0.273
$ws.Range("AN3").Value = 5.172413793103448
$ws.Range("AO3").Value = -0.8763636363636362
$ws.Range("AP3").Value = 5.027586534975698
$ws.Range("AQ3").Value = -0.8827838827838826
$ws.Range("A4").Value = "Netherlands"
$ws.Range("B4").Value = "Fastned B.V. (ENXTAM:FAST)"
$ws.Range("C4").Value = "Power"
$ws.Range("G4").Value = -0.6942590120160214
$ws.Range("H4").Value = -0.6942590120160214
$ws.Range("I4").Value = -1.228304405874499
$ws.Range("J4").Value = -1.228304405874499
$ws.Range("K4").Value = -13.8
$ws.Range("L4").Value = -1.842456608811749
$ws.Range("M4").Value = -0
$ws.Range("N4").Value = -0
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = -0
$ws.Range("Q4").Value = -0
$ws.Range("R4").Value = 0
$ws.Range("S4").Value = 0
$ws.Range("U4").Value = 15.2
$ws.Range("V4").Value = 0.01805011281320508
$ws.Range("W4").Value = 5.227272727272728
$ws.Range("X4").Value = 0.03484300448517595
$ws.Range("Y4").Value = 5.192429722787551
$ws.Range("Z4").Value = 0.2154775604142693
$ws.Range("AA4").Value = -0.2646720368239355
$ws.Range("AB4").Value = 0.03367754126616151
$ws.Range("AC4").Value = -0.298349578090097
$ws.Range("AD4").Value = 68.5
$ws.Range("AE4").Value = 0
$ws.Range("AF4").Value = 68.5
$ws.Range("AG4").Value = 53.3
$ws.Range("AH4").Value = 0.075225126290358
$ws.Range("AI4").Value = 1.273234200743494
$ws.Range("AJ4").Value = 0.0595264686173777
$ws.Range("AK4").Value = 1.380829015544041
$ws.Range("AL4").Value = 4
$ws.Range("AM4").Value = 3.868
$ws.Range("AN4").Value = -11.79001721170396
$ws.Range("AO4").Value = -2.3
$ws.Range("AP4").Value = -9.173838209982788
$ws.Range("AQ4").Value = -2.378490175801448
